$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 4th voltage-scale part added to the BOM: REFERENCE VOLTAGE / Microchip / qty 1
$ws.Range("A11").Value = "REFERENCE VOLTAGE"
$ws.Range("B11").Value = "Microchip"
$ws.Range("C11").Value = 1

# Link the new part to its datasheet/product page
$ws.Hyperlinks.Add($ws.Range("A11"), "https://www.microchip.com/en-us/product/MCP1525")

# Give A11 the same "Hyperlink" look as the other linked parts (A3:A10)
$ws.Range("A11").Style = "Hyperlink"

# Column A's text got longer ("REFERENCE VOLTAGE"/"SHIFT REGISTER" no longer the
# longest) so it was manually widened instead of left on AutoFit/BestFit.
$ws.Columns.Item(1).ColumnWidth = 22

# Leave the selection where the user's cursor ended up after entering the row
$ws.Range("B12").Select() | Out-Null
